# Auto-generated edit script: add rows 206-220 and update dimension/col width
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A206").Value = 'Why can''t I add 251 curve shades to my log?'
$ws.Range("B206").Value = 'llama3.2:latest'
$ws.Range("C206").Value = 'You cannot add 251 curve shades to your log because the maximum number of curve shades per plot is 250.'

$ws.Range("A207").Value = 'What''s the maximum number of data points allowed per curve?'
$ws.Range("B207").Value = 'llama3.2:latest'
$ws.Range("C207").Value = 'Unfortunately, I was unable to find any information on the maximum number of data points allowed per curve in the provided documentation.'

$ws.Range("A208").Value = 'What''s the maximum number of data points allowed per curve?'
$ws.Range("B208").Value = 'llama3.2:latest'
$ws.Range("C208").Value = 'Unfortunately, I was unable to find any information on the maximum number of data points allowed per curve in the provided documentation.'

$ws.Range("A209").Value = 'What''s the maximum curves I can load in a data file?'
$ws.Range("B209").Value = 'llama3.2:latest'
$ws.Range("C209").Value = 'According to the Document section, the answer is:
You can load an unlimited number of data files.'

$ws.Range("A210").Value = 'What''s the maximum number of data points allowed per curve?'
$ws.Range("B210").Value = 'llama3.2:latest'
$ws.Range("C210").Value = 'Unfortunately, I was unable to find any information on the maximum number of data points allowed per curve in the provided documentation.'

$ws.Range("A211").Value = 'What''s the maximum number of data points allowed per curve?'
$ws.Range("B211").Value = 'llama3.2:latest'
$ws.Range("C211").Value = 'The maximum number of data points allowed per curve is unlimited.'

$ws.Range("A212").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B212").Value = 'llama3.2:latest'
$ws.Range("C212").Value = 'To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Click "OK" to apply the changes.'

$ws.Range("A213").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B213").Value = 'llama3.2:latest'
$ws.Range("C213").Value = 'To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Click "OK" to apply the changes.'

$ws.Range("A214").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B214").Value = 'llama3.2:latest'
$ws.Range("C214").Value = 'To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Click "OK" to apply the changes.'

$ws.Range("A215").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B215").Value = 'llama3.2:latest'
$ws.Range("C215").Value = 'To set the curve shading name, click on the "Curve" tab in the GEO application and select "Shading" from the dropdown menu. Then, enter "Hydrocarbon bearing zone highlighted" in the "Shading Name" field and press Enter.'

$ws.Range("A216").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B216").Value = 'llama3.2:latest'
$ws.Range("C216").Value = 'To set the curve shading name, click on the "Curve" tab in the GEO application and select "Shading" from the dropdown menu. Then, enter "Hydrocarbon bearing zone highlighted" in the "Shading Name" field and press Enter.'

$ws.Range("A217").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B217").Value = 'llama3.2:latest'
$ws.Range("C217").Value = 'To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Click "OK" to apply the changes.
Note: The character length of the curve shade name should be within the allowed limit to avoid any errors or issues with the application.'

$ws.Range("A218").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B218").Value = 'llama3.2:latest'
$ws.Range("C218").Value = 'To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Click "OK" to apply the changes.
Note: The character length of the curve shade name should be within the allowed limit to avoid any errors or issues with the application.'

$ws.Range("A219").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B219").Value = 'llama3.2:latest'
$ws.Range("C219").Value = 'Based on the provided feedback, here''s an optimized answer:
**I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"**
To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field. Note that the character limit for the curve shade name is 20 characters; ensure your chosen name does not exceed this limit to avoid any errors or issues with the application.
This answer addresses the feedback by:
1. Providing a clear step-by-step guide on how to set the curve shading name.
2. Mentioning the character limit of 20 characters for the curve shade name, which was previously mentioned in one of the feedback responses.
3. Encouraging the user to check if their chosen name exceeds the character limit and make a decision accordingly.
This optimized answer should provide a more accurate and helpful response to the user''s question while also addressing the concerns raised in the feedback.'

$ws.Range("A220").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B220").Value = 'llama3.2:latest'
$ws.Range("C220").Value = 'Here is an optimized answer:
I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted".
To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field.
Note that the character limit for the curve shade name is 20 characters. Since your chosen name has a length of 37 characters, it exceeds the allowed limit. Therefore, this operation is not allowed due to the character length constraint.
This answer addresses the feedback by:
1. Providing a clear step-by-step guide on how to set the curve shading name.
2. Mentioning the character limit of 20 characters for the curve shade name and explaining why the operation is not allowed due to exceeding this limit.
3. Encouraging the user to check if their chosen name exceeds the character limit and make a decision accordingly.
This optimized answer should provide a more accurate and helpful response to the user''s question while also addressing the concerns raised in the feedback.'

$ws.Range("A206:A220").EntireRow.AutoFit()

$ws.Columns.Item(3).ColumnWidth = 1089 - (5/6)
